# Pet Rules.xlsx edit
# - G7: $p.setModifiedRatingFactor($param) -> $p.setModifiedRatingFactor("$param");  (curly quotes)
# - H7: $p.setTier($param) -> $p.setTier("$param");  (curly quotes)
# - Selection moved from A5 to F29
# - Column G width widened (30.96 -> 33.39 "Calc" units)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$leftCurly = [char]0x201C
$rightCurly = [char]0x201D

$ws.Range("G7").Value = '$p.setModifiedRatingFactor(' + $leftCurly + '$param' + $rightCurly + ');'
$ws.Range("H7").Value = '$p.setTier(' + $leftCurly + '$param' + $rightCurly + ');'

# Widen column G (closest achievable stored width to the target 33.39)
$ws.Columns.Item(7).ColumnWidth = 32.5

# Move the active cell/selection to F29
$ws.Range("F29").Select()
